$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 144 (existing rows 144:173 shift down to 145:174)
$ws.Rows.Item(144).Insert()

# Populate the newly inserted row 144 with a new weekly price observation,
# mirroring the surrounding record (row 143) except for the date and volume.
$ws.Cells.Item(144, 1).Value = 10
$ws.Cells.Item(144, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(144, 3).Value = "La Araucanía"
$ws.Cells.Item(144, 4).Value = 44511
$ws.Cells.Item(144, 5).Value = 9
$ws.Cells.Item(144, 6).Value = 100112039
$ws.Cells.Item(144, 7).Value = "Ciboulette"
$ws.Cells.Item(144, 8).Value = "Sin especificar"
$ws.Cells.Item(144, 9).Value = "Primera"
$ws.Cells.Item(144, 10).Value = 50
$ws.Cells.Item(144, 11).Value = 5000
$ws.Cells.Item(144, 12).Value = 5000
$ws.Cells.Item(144, 13).Value = 5000
$ws.Cells.Item(144, 14).Value = "$/docena de atados"
$ws.Cells.Item(144, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(144, 16).Value = 1667
$ws.Cells.Item(144, 17).Value = 3
$ws.Cells.Item(144, 18).Value = "Hortaliza"
